$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Constants (so we don't depend on the PowerShell-only Xl* enum names)
# ---------------------------------------------------------------------------
$xlContinuous   = 1
$xlMedium       = -4138
$xlCenter       = -4108
$xlTop          = -4160
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1. Fill in the new "UNITS" column (G) for both tables first, then the QCP /
#    HumMod column headers (H1/H10), then the footnote. This mirrors the
#    order new shared strings were introduced in the target workbook
#    (UNITS, mmHg, L/min, QCP, HumMod, note).
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "UNITS"
$ws.Range("G3").Value = "mmHg"
$ws.Range("G5").Value = "mmHg"
$ws.Range("G7").Value = "L/min"

$ws.Range("G10").Value = "UNITS"
$ws.Range("G12").Value = "mmHg"
$ws.Range("G14").Value = "mmHg"
$ws.Range("G16").Value = "L/min"

$ws.Range("H1").Value = "QCP"
$ws.Range("H10").Value = "HumMod"

# ---------------------------------------------------------------------------
# 2. Second table (HumMod data), rows 10-16, mirrors rows 1-7
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Time (Min)"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 20

$ws.Range("A11").Value = "Venous pH"
$ws.Range("B11").Value = 7.38
$ws.Range("C11").Value = 7.55
$ws.Range("D11").Value = 7.65
$ws.Range("E11").Value = 7.73
$ws.Range("F11").Value = 7.79

$ws.Range("A12").Value = "Venous pCO2"
$ws.Range("B12").Value = 41.6
$ws.Range("C12").Value = 39.3
$ws.Range("D12").Value = 39.1
$ws.Range("E12").Value = 38.6
$ws.Range("F12").Value = 37.9

$ws.Range("A13").Value = "Arterial pH"
$ws.Range("B13").Value = 7.43
$ws.Range("C13").Value = 7.58
$ws.Range("D13").Value = 7.68
$ws.Range("E13").Value = 7.75
$ws.Range("F13").Value = 7.81

$ws.Range("A14").Value = "Arterial pCO2"
$ws.Range("B14").Value = 37.4
$ws.Range("C14").Value = 36.4
$ws.Range("D14").Value = 36.7
$ws.Range("E14").Value = 36.5
$ws.Range("F14").Value = 36

$ws.Range("A15").Value = "Brain pH"
$ws.Range("B15").Value = 7.12
$ws.Range("C15").Value = 7.16
$ws.Range("D15").Value = 7.17
$ws.Range("E15").Value = 7.19
$ws.Range("F15").Value = 7.21

$ws.Range("A16").Value = "Ventilation"
$ws.Range("B16").Value = 6.6
$ws.Range("C16").Value = 4.2
$ws.Range("D16").Value = 3.6
$ws.Range("E16").Value = 3.5
$ws.Range("F16").Value = 3.3

# ---------------------------------------------------------------------------
# 3. Note paragraph under the tables (merged A18:G19, yellow fill)
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "For the HumMod version, the IV settings were an H2O rate of 100 mL/min with a bicarbonate concentration of 500 mmol/L"

# ---------------------------------------------------------------------------
# 4. Build "master" cell formats from scratch on blank cells (so the style
#    table stays compact and every identically formatted cell shares one
#    style index), then copy those formats onto all the relevant ranges.
# ---------------------------------------------------------------------------

# 4a. "Row label" style (column A): Arial 12, wrap, vertical-top, full box border
$masterA = $ws.Range("J30")
$masterA.Font.Name = "Arial"
$masterA.Font.Size = 12
$masterA.WrapText = $true
$masterA.VerticalAlignment = $xlTop
$masterA.Borders.LineStyle = $xlContinuous
$masterA.Borders.Weight = $xlMedium

# 4b. "Data value" style (columns B:F): same as above + centered horizontally
$masterB = $ws.Range("J31")
$masterB.Font.Name = "Arial"
$masterB.Font.Size = 12
$masterB.WrapText = $true
$masterB.VerticalAlignment = $xlTop
$masterB.HorizontalAlignment = $xlCenter
$masterB.Borders.LineStyle = $xlContinuous
$masterB.Borders.Weight = $xlMedium

# 4c. "Units" style (column G): default font, default alignment, full box border
$masterG = $ws.Range("J32")
$masterG.Borders.LineStyle = $xlContinuous
$masterG.Borders.Weight = $xlMedium

# 4d. Note-heading style (A18): Arial 12, wrap, vertical-top, centered, yellow fill
$masterNote = $ws.Range("J33")
$masterNote.Font.Name = "Arial"
$masterNote.Font.Size = 12
$masterNote.WrapText = $true
$masterNote.VerticalAlignment = $xlTop
$masterNote.HorizontalAlignment = $xlCenter
$masterNote.Interior.Color = 65535

# 4e. Note-fill style (B18:G19/A19:G19): default font, centered, yellow fill
$masterNoteFill = $ws.Range("J34")
$masterNoteFill.HorizontalAlignment = $xlCenter
$masterNoteFill.Interior.Color = 65535

# ---------------------------------------------------------------------------
# 5. Propagate the master formats onto the real ranges
# ---------------------------------------------------------------------------
$masterA.Copy()
$ws.Range("A1:A7").PasteSpecial($xlPasteFormats)
$ws.Range("A10:A16").PasteSpecial($xlPasteFormats)

$masterB.Copy()
$ws.Range("B1:F7").PasteSpecial($xlPasteFormats)
$ws.Range("B10:F16").PasteSpecial($xlPasteFormats)

$masterG.Copy()
$ws.Range("G1:G7").PasteSpecial($xlPasteFormats)
$ws.Range("G10:G16").PasteSpecial($xlPasteFormats)

$masterNote.Copy()
$ws.Range("A18").PasteSpecial($xlPasteFormats)

$masterNoteFill.Copy()
$ws.Range("B18:G19").PasteSpecial($xlPasteFormats)
$ws.Range("A19").PasteSpecial($xlPasteFormats)

# Clear out the scratch master cells so they don't leave stray formatting
$ws.Range("J30:J34").Clear()

# ---------------------------------------------------------------------------
# 6. Row heights (match the existing table's look) + merge + blank spacer row
# ---------------------------------------------------------------------------
$ws.Rows.Item(9).RowHeight = 15.75

$ws.Rows.Item(10).RowHeight = 15.75
$ws.Rows.Item(11).RowHeight = 15.75
$ws.Rows.Item(12).RowHeight = 30.75
$ws.Rows.Item(13).RowHeight = 15.75
$ws.Rows.Item(14).RowHeight = 30.75
$ws.Rows.Item(15).RowHeight = 15.75
$ws.Rows.Item(16).RowHeight = 15.75

$ws.Range("A18:G19").Merge()

# ---------------------------------------------------------------------------
# 7. Selection, matching the diff's recorded sheet view state
# ---------------------------------------------------------------------------
$ws.Range("G12").Select()
